$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows per the diff
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -2
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = -3
$ws.Range("F14").Value = -4
$ws.Range("F19").Value = 1
$ws.Range("F22").Value = -1
$ws.Range("F28").Value = -2
